# Refined metadata to be additional tab
#
# 1. Update the "panel_query_time" (column F) timestamps on the "data" sheet
#    to the values recorded for the later query run.
# 2. Add a new "metadata" worksheet (after "data") holding one summary row
#    describing the panel query itself.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the per-row query timestamps on the "data" sheet ----------
$newTimes = @(
    "2021-10-05 14:35:08.611804",
    "2021-10-05 14:35:08.611812",
    "2021-10-05 14:35:08.611815",
    "2021-10-05 14:35:08.611818",
    "2021-10-05 14:35:08.611821",
    "2021-10-05 14:35:08.611824",
    "2021-10-05 14:35:08.611826",
    "2021-10-05 14:35:08.611829",
    "2021-10-05 14:35:08.611832",
    "2021-10-05 14:35:08.611835",
    "2021-10-05 14:35:08.611837",
    "2021-10-05 14:35:08.611840",
    "2021-10-05 14:35:08.611842",
    "2021-10-05 14:35:08.611845",
    "2021-10-05 14:35:08.611847",
    "2021-10-05 14:35:08.611850",
    "2021-10-05 14:35:08.611853",
    "2021-10-05 14:35:08.611856",
    "2021-10-05 14:35:08.611859",
    "2021-10-05 14:35:08.611862",
    "2021-10-05 14:35:08.611865",
    "2021-10-05 14:35:08.611868",
    "2021-10-05 14:35:08.611870",
    "2021-10-05 14:35:08.611873",
    "2021-10-05 14:35:08.611876",
    "2021-10-05 14:35:08.611878",
    "2021-10-05 14:35:08.611881",
    "2021-10-05 14:35:08.611884",
    "2021-10-05 14:35:08.611886",
    "2021-10-05 14:35:08.611889"
)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
    $dataSheet.Cells.Item($i + 2, 6).Value = $newTimes[$i]
}

# --- 2. Add the "metadata" sheet, placed right after "data" ---------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "metadata"
# Look up "data" fresh (a reference captured earlier can go stale and
# silently no-op the Move), then reposition "metadata" to follow it.
$newSheet.Move($null, $wb.Worksheets.Item("data"))

# Re-resolve "metadata" by name after the Move: sheet handles captured
# before a sheet-order mutation (Add/Move) can silently re-target a
# different sheet afterwards, so every write below uses a fresh lookup.
$metaSheet = $wb.Worksheets.Item("metadata")

# Header row (bold, centered, top-aligned, thin border) - mirrors the style
# used for the header row on the "data" sheet.
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $cell = $metaSheet.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Row-index cell (A2), styled the same way as the "data" sheet's index column.
$indexCell = $metaSheet.Cells.Item(2, 1)
$indexCell.Value = 0
$indexCell.Font.Bold = $true
$indexCell.HorizontalAlignment = -4108
$indexCell.VerticalAlignment = -4160
$indexCell.Borders.LineStyle = 1

# Data row describing the panel query.
$metaSheet.Cells.Item(2, 2).Value = "Overgrowth"
$metaSheet.Cells.Item(2, 3).Value = 151
# "data_version" is a text value ("1.4") in the source data, not a number -
# force Text format first so the COM layer doesn't silently coerce it.
$versionCell = $metaSheet.Cells.Item(2, 4)
$versionCell.NumberFormat = "@"
$versionCell.Value = "1.4"
$metaSheet.Cells.Item(2, 5).Value = "2021-08-31T01:05:58.736855Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:35:08.608203"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/151/?format=json"

Write-Output "metadata sheet added; data timestamps refreshed"
